$d = $word.ActiveDocument

# --- 1. Refresh the footer timestamp ------------------------------------
foreach ($sec in $d.Sections) {
    $footer = $sec.Footers(1)
    $footer.Range.Find.Execute(
        "2025-06-30 12:12Z / ", $true, $false, $false, $false, $false,
        $true, 1, $false, "2025-07-02 02:48Z / ", 2
    )
}

# --- 2. Add the regression-test character styles (b, i, sub, sup, u) ----
$wdStyleTypeCharacter = 2
$wdUnderlineSingle = 1

function Add-CharStyle([string]$name) {
    $style = $d.Styles.Add($name, $wdStyleTypeCharacter)
    $style.BaseStyle = "DefaultParagraphFont"
    $style.Priority = 1
    $style.QuickStyle = $true
    return $style
}

$b = Add-CharStyle "b"
$b.Font.Bold = $true

$i = Add-CharStyle "i"
$i.Font.Italic = $true

$sub = Add-CharStyle "sub"
$sub.Font.Subscript = $true

$sup = Add-CharStyle "sup"
$sup.Font.Superscript = $true

$u = Add-CharStyle "u"
$u.Font.Underline = $wdUnderlineSingle
